# #327 Ajout des profils d'acces a58d18c1e8091c98efec92c8c093b361a253eee5
#
# 1. Metadata sheet: bump the generation "Date" value (B8) to the new
#    timestamp.
# 2. Elements sheet: the "Mapping" columns were reordered so that the
#    "Spécification métier vers l'extension ROR ReopeningDate" mapping
#    now comes before the "RIM Mapping" one -- i.e. columns AK (37) and
#    AL (38) are swapped in their entirety (width, header text, and every
#    data cell in the used range).

$wb = $excel.ActiveWorkbook

# --- 1. Metadata!B8 date bump ------------------------------------------------
$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B8").Value2 = "2024-03-19T13:17:15+00:00"

# --- 2. Elements: swap columns AK (37) and AL (38) --------------------------
$elements = $wb.Worksheets.Item("Elements")

# Column widths (values taken from the original sheet's <cols> entries,
# simply swapped between the two columns). The runtime's ColumnWidth
# setter expects "characters", which is the stored XML width minus the
# standard ~0.8333 padding.
$padding = 0.8333333333333333
$akOriginalWidth = 24.98046875
$alOriginalWidth = 70.94140625

$elements.Columns.Item(37).ColumnWidth = $alOriginalWidth - $padding
$elements.Columns.Item(38).ColumnWidth = $akOriginalWidth - $padding

# Header + every data row in the used range (rows 1-6) -- swap the two
# columns' contents cell by cell. Skip rows where both cells already hold
# the same value (e.g. both blank) so untouched rows stay byte-identical
# instead of being rewritten as "empty" cells.
$lastRow = $elements.UsedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $elements.Cells.Item($r, 37)
    $alCell = $elements.Cells.Item($r, 38)
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    if ($akVal -ne $alVal) {
        $akCell.Value2 = $alVal
        $alCell.Value2 = $akVal
    }
}
